$d = $word.ActiveDocument

$replacements = @(
    @("704×7=4928", "170×7=1190"),
    @("867×2=1734", "456×6=2736"),
    @("790×3=2370", "629×3=1887"),
    @("678×7=4746", "388×3=1164"),
    @("511×2=1022", "880×6=5280"),
    @("525×7=3675", "680×6=4080"),
    @("862×8=6896", "893×8=7144"),
    @("838×9=7542", "519×2=1038"),
    @("339×9=3051", "636×6=3816"),
    @("819×3=2457", "794×2=1588"),
    @("727×6=4362", "346×4=1384"),
    @("505×7=3535", "137×6=822"),
    @("144×2=288", "289×9=2601"),
    @("655×7=4585", "347×7=2429"),
    @("147×9=1323", "190×7=1330"),
    @("182×9=1638", "838×7=5866"),
    @("733×4=2932", "867×9=7803"),
    @("413×9=3717", "949×4=3796"),
    @("359×9=3231", "686×6=4116"),
    @("894×6=5364", "306×5=1530"),
    @("993×6=5958", "752×9=6768"),
    @("329×7=2303", "915×8=7320"),
    @("283×4=1132", "698×4=2792"),
    @("786×2=1572", "626×7=4382"),
    @("408×8=3264", "848×9=7632")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
